$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-11-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-14 Thursday", 2) | Out-Null

# Update each table cell by explicit row/column to avoid ambiguity from duplicate values
$t = $d.Tables(1)
$t.Cell(1, 1).Range.Text = "77-52=25"
$t.Cell(1, 2).Range.Text = "7+40=47"
$t.Cell(1, 3).Range.Text = "12+48=60"
$t.Cell(1, 4).Range.Text = "88-79=9"
$t.Cell(1, 5).Range.Text = "70+10=80"
$t.Cell(2, 1).Range.Text = "81-51=30"
$t.Cell(2, 2).Range.Text = "0+70=70"
$t.Cell(2, 3).Range.Text = "23+18=41"
$t.Cell(2, 4).Range.Text = "40+8=48"
$t.Cell(2, 5).Range.Text = "55-35=20"
$t.Cell(3, 1).Range.Text = "30+0=30"
$t.Cell(3, 2).Range.Text = "71-49=22"
$t.Cell(3, 3).Range.Text = "95-69=26"
$t.Cell(3, 4).Range.Text = "55-22=33"
$t.Cell(3, 5).Range.Text = "60-2=58"
$t.Cell(4, 1).Range.Text = "22+22=44"
$t.Cell(4, 2).Range.Text = "70-25=45"
$t.Cell(4, 3).Range.Text = "74-29=45"
$t.Cell(4, 4).Range.Text = "51+19=70"
$t.Cell(4, 5).Range.Text = "32-19=13"
$t.Cell(5, 1).Range.Text = "28-21=7"
$t.Cell(5, 2).Range.Text = "3+50=53"
$t.Cell(5, 3).Range.Text = "28+58=86"
$t.Cell(5, 4).Range.Text = "21+19=40"
$t.Cell(5, 5).Range.Text = "8+68=76"
$t.Cell(6, 1).Range.Text = "74-55=19"
$t.Cell(6, 2).Range.Text = "0+80=80"
$t.Cell(6, 3).Range.Text = "40+52=92"
$t.Cell(6, 4).Range.Text = "40+10=50"
$t.Cell(6, 5).Range.Text = "47-31=16"
$t.Cell(7, 1).Range.Text = "64-59=5"
$t.Cell(7, 2).Range.Text = "31+56=87"
$t.Cell(7, 3).Range.Text = "46-15=31"
$t.Cell(7, 4).Range.Text = "43-14=29"
$t.Cell(7, 5).Range.Text = "55-9=46"
$t.Cell(8, 1).Range.Text = "35+18=53"
$t.Cell(8, 2).Range.Text = "34+0=34"
$t.Cell(8, 3).Range.Text = "43-27=16"
$t.Cell(8, 4).Range.Text = "37+26=63"
$t.Cell(8, 5).Range.Text = "12+7=19"
$t.Cell(9, 1).Range.Text = "92-0=92"
$t.Cell(9, 2).Range.Text = "4+47=51"
$t.Cell(9, 3).Range.Text = "7+48=55"
$t.Cell(9, 4).Range.Text = "30-0=30"
$t.Cell(9, 5).Range.Text = "95-15=80"
$t.Cell(10, 1).Range.Text = "4+74=78"
$t.Cell(10, 2).Range.Text = "52+33=85"
$t.Cell(10, 3).Range.Text = "98-27=71"
$t.Cell(10, 4).Range.Text = "84+11=95"
$t.Cell(10, 5).Range.Text = "75-4=71"
$t.Cell(11, 1).Range.Text = "96-8=88"
$t.Cell(11, 2).Range.Text = "77-74=3"
$t.Cell(11, 3).Range.Text = "86-76=10"
$t.Cell(11, 4).Range.Text = "39+6=45"
$t.Cell(11, 5).Range.Text = "60+1=61"
$t.Cell(12, 1).Range.Text = "24+72=96"
$t.Cell(12, 2).Range.Text = "72-41=31"
$t.Cell(12, 3).Range.Text = "32+61=93"
$t.Cell(12, 4).Range.Text = "20+45=65"
$t.Cell(12, 5).Range.Text = "48+27=75"
$t.Cell(13, 1).Range.Text = "28+22=50"
$t.Cell(13, 2).Range.Text = "71+22=93"
$t.Cell(13, 3).Range.Text = "74-20=54"
$t.Cell(13, 4).Range.Text = "39+35=74"
$t.Cell(13, 5).Range.Text = "43-5=38"
$t.Cell(14, 1).Range.Text = "96-42=54"
$t.Cell(14, 2).Range.Text = "32+41=73"
$t.Cell(14, 3).Range.Text = "58-44=14"
$t.Cell(14, 4).Range.Text = "77-7=70"
$t.Cell(14, 5).Range.Text = "35+6=41"
$t.Cell(15, 1).Range.Text = "85-20=65"
$t.Cell(15, 2).Range.Text = "32+6=38"
$t.Cell(15, 3).Range.Text = "39-25=14"
$t.Cell(15, 4).Range.Text = "34-5=29"
$t.Cell(15, 5).Range.Text = "99-5=94"
$t.Cell(16, 1).Range.Text = "41+44=85"
$t.Cell(16, 2).Range.Text = "51-49=2"
$t.Cell(16, 3).Range.Text = "91-18=73"
$t.Cell(16, 4).Range.Text = "28-5=23"
$t.Cell(16, 5).Range.Text = "14+37=51"
$t.Cell(17, 1).Range.Text = "61+20=81"
$t.Cell(17, 2).Range.Text = "30+62=92"
$t.Cell(17, 3).Range.Text = "90-72=18"
$t.Cell(17, 4).Range.Text = "70-55=15"
$t.Cell(17, 5).Range.Text = "22+32=54"
$t.Cell(18, 1).Range.Text = "28+49=77"
$t.Cell(18, 2).Range.Text = "76+19=95"
$t.Cell(18, 3).Range.Text = "44-42=2"
$t.Cell(18, 4).Range.Text = "7+45=52"
$t.Cell(18, 5).Range.Text = "95-88=7"
$t.Cell(19, 1).Range.Text = "20-13=7"
$t.Cell(19, 2).Range.Text = "33+42=75"
$t.Cell(19, 3).Range.Text = "59+14=73"
$t.Cell(19, 4).Range.Text = "56-9=47"
$t.Cell(19, 5).Range.Text = "99-85=14"
$t.Cell(20, 1).Range.Text = "60-2=58"
$t.Cell(20, 2).Range.Text = "88-74=14"
$t.Cell(20, 3).Range.Text = "40+37=77"
$t.Cell(20, 4).Range.Text = "87-60=27"
$t.Cell(20, 5).Range.Text = "13-10=3"
